$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 261
$dates = @(
    "17-09-2021",
    "18-09-2021",
    "19-09-2021",
    "20-09-2021",
    "21-09-2021",
    "22-09-2021",
    "23-09-2021",
    "24-09-2021",
    "25-09-2021",
    "26-09-2021",
    "27-09-2021",
    "28-09-2021",
    "29-09-2021",
    "30-09-2021"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = 3623
    $ws.Cells.Item($r, 3).Value = 240
}
